$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7342893.5
$ws.Range("J17").Value = 7342893.5
$ws.Range("L17").Value = 22028680.5
$ws.Range("N17").Value = -22029016.5
$ws.Range("H52").Value = 999
$ws.Range("I52").Value = 999
$ws.Range("K52").Value = 2997
$ws.Range("M52").Value = -2837
$ws.Range("H53").Value = 712.9583
$ws.Range("I53").Value = 666.35297
$ws.Range("J53").Value = 826.1429000000001
$ws.Range("K53").Value = 666.35297
$ws.Range("L53").Value = 826.1429000000001
$ws.Range("M53").Value = -29.35297000000003
$ws.Range("N53").Value = -2100.1429
$ws.Range("H98").Value = 1799
$ws.Range("I98").Value = 1623.75
$ws.Range("K98").Value = 1623.75
$ws.Range("M98").Value = -125.75
$ws.Range("H122").Value = 1799
$ws.Range("I122").Value = 1623.75
$ws.Range("K122").Value = 4871.25
$ws.Range("M122").Value = -2421.25
$ws.Range("H137").Value = 5560.3447
$ws.Range("I137").Value = 5559.5654
$ws.Range("K137").Value = 16678.6962
$ws.Range("M137").Value = -14128.6962
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 49.666668
$ws.Range("I5").Value = 49.666668
$ws.Range("K5").Value = 49.666668
$ws.Range("M5").Value = 62.333332
$ws.Range("H46").Value = 10470.2
$ws.Range("J46").Value = 10470.2
$ws.Range("L46").Value = 10470.2
$ws.Range("N46").Value = -11108.2
$ws.Range("H50").Value = 2048.3333
$ws.Range("J50").Value = 2048.3333
$ws.Range("L50").Value = 2048.3333
$ws.Range("N50").Value = -3476.3333
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 49.666668
$ws.Range("I4").Value = 49.666668
$ws.Range("K4").Value = 49.666668
$ws.Range("M4").Value = 65.333332
$ws.Range("H8").Value = 1614.6666
$ws.Range("J8").Value = 1900
$ws.Range("L8").Value = 1900
$ws.Range("N8").Value = -2180
$ws.Range("H96").Value = 35749.5
$ws.Range("I96").Value = 35749.5
$ws.Range("K96").Value = 35749.5
$ws.Range("M96").Value = -33003.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 26057.715
$ws.Range("I69").Value = 17480.8
$ws.Range("J69").Value = 47500
$ws.Range("K69").Value = 17480.8
$ws.Range("L69").Value = 47500
$ws.Range("M69").Value = -16731.8
$ws.Range("N69").Value = -48998
$ws.Range("H72").Value = 26057.715
$ws.Range("I72").Value = 17480.8
$ws.Range("J72").Value = 47500
$ws.Range("K72").Value = 52442.39999999999
$ws.Range("L72").Value = 142500
$ws.Range("M72").Value = -48698.39999999999
$ws.Range("N72").Value = -149988
$ws.Range("H74").Value = 42739.25
$ws.Range("J74").Value = 42739.25
$ws.Range("L74").Value = 42739.25
$ws.Range("N74").Value = -44487.25
$ws.Range("H77").Value = 42739.25
$ws.Range("J77").Value = 42739.25
$ws.Range("L77").Value = 128217.75
$ws.Range("N77").Value = -136953.75
$ws.Range("H98").Value = 89900
$ws.Range("J98").Value = 89900
$ws.Range("L98").Value = 89900
$ws.Range("N98").Value = -94392
$ws.Range("H105").Value = 1682.2916
$ws.Range("I105").Value = 1308.4286
$ws.Range("K105").Value = 1308.4286
$ws.Range("M105").Value = 438.5714
$ws.Range("H110").Value = 79875
$ws.Range("J110").Value = 79875
$ws.Range("L110").Value = 79875
$ws.Range("N110").Value = -88055
$ws.Range("H111").Value = 37000
$ws.Range("J111").Value = 37000
$ws.Range("L111").Value = 37000
$ws.Range("N111").Value = -45180
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1800
$ws.Range("I17").Value = 200
$ws.Range("K17").Value = 600
$ws.Range("M17").Value = -431
$ws.Range("H39").Value = 1550
$ws.Range("I39").Value = 1200
$ws.Range("J39").Value = 1666.6666
$ws.Range("K39").Value = 3600
$ws.Range("L39").Value = 4999.9998
$ws.Range("M39").Value = -3306
$ws.Range("N39").Value = -5587.9998
$ws.Range("H55").Value = 6936.8335
$ws.Range("J55").Value = 7476.636
$ws.Range("L55").Value = 22429.908
$ws.Range("N55").Value = -22783.908
$ws.Range("H107").Value = 1272.0714
$ws.Range("J107").Value = 1400.9
$ws.Range("L107").Value = 4202.700000000001
$ws.Range("N107").Value = -8042.700000000001
$ws.Range("H110").Value = 2727
$ws.Range("I110").Value = 2727
$ws.Range("K110").Value = 8181
$ws.Range("M110").Value = -4091
$ws.Range("I117").Value = 47624976
$ws.Range("J117").Value = 912751.4399999999
$ws.Range("K117").Value = 142874928
$ws.Range("L117").Value = 2738254.32
$ws.Range("M117").Value = -142871486
$ws.Range("N117").Value = -2745138.32
$ws.Range("H122").Value = 5657.543
$ws.Range("J122").Value = 6758.36
$ws.Range("L122").Value = 60825.24
$ws.Range("N122").Value = -65725.23999999999
$ws.Range("H141").Value = 9409.799999999999
$ws.Range("I141").Value = 9409.799999999999
$ws.Range("K141").Value = 28229.4
$ws.Range("M141").Value = -23049.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 57999.668
$ws.Range("J51").Value = 57599.8
$ws.Range("L51").Value = 57599.8
$ws.Range("N51").Value = -58617.8
$ws.Range("H80").Value = 3088.5557
$ws.Range("I80").Value = 1949.5
$ws.Range("J80").Value = 3999.8
$ws.Range("K80").Value = 1949.5
$ws.Range("L80").Value = 3999.8
$ws.Range("M80").Value = -951.5
$ws.Range("N80").Value = -5995.8
$ws.Range("H83").Value = 3088.5557
$ws.Range("I83").Value = 1949.5
$ws.Range("J83").Value = 3999.8
$ws.Range("K83").Value = 9747.5
$ws.Range("L83").Value = 19999
$ws.Range("M83").Value = -4755.5
$ws.Range("N83").Value = -29983
$ws.Range("H122").Value = 5103
$ws.Range("I122").Value = 5103
$ws.Range("K122").Value = 15309
$ws.Range("M122").Value = -12859
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 25095000
$ws.Range("J20").Value = 25095000
$ws.Range("L20").Value = 25095000
$ws.Range("N20").Value = -25095452
$ws.Range("H31").Value = 4516.125
$ws.Range("J31").Value = 4399.75
$ws.Range("L31").Value = 4399.75
$ws.Range("N31").Value = -4895.75
$ws.Range("H82").Value = 144754.42
$ws.Range("I82").Value = 2176.2
$ws.Range("J82").Value = 501200
$ws.Range("K82").Value = 2176.2
$ws.Range("L82").Value = 501200
$ws.Range("M82").Value = -1815.2
$ws.Range("N82").Value = -501922
$ws.Range("H85").Value = 144754.42
$ws.Range("I85").Value = 2176.2
$ws.Range("J85").Value = 501200
$ws.Range("K85").Value = 2176.2
$ws.Range("L85").Value = 501200
$ws.Range("M85").Value = -928.1999999999998
$ws.Range("N85").Value = -503696
$ws.Range("H122").Value = 3865.276
$ws.Range("I122").Value = 3904.2917
$ws.Range("K122").Value = 11712.8751
$ws.Range("M122").Value = -9262.875100000001
$ws.Range("H132").Value = 10593.478
$ws.Range("I132").Value = 11745.875
$ws.Range("J132").Value = 7520.4165
$ws.Range("K132").Value = 35237.625
$ws.Range("L132").Value = 22561.2495
$ws.Range("M132").Value = -32707.625
$ws.Range("N132").Value = -27621.2495
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").ClearContents()
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = 0
$ws.Range("H28").Value = 14714.25
$ws.Range("I28").Value = 13950
$ws.Range("J28").Value = 14969
$ws.Range("K28").Value = 13950
$ws.Range("L28").Value = 14969
$ws.Range("M28").Value = -13602
$ws.Range("N28").Value = -15665
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").ClearContents()
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = 0
$ws.Range("H107").Value = 2349.9412
$ws.Range("I107").Value = 1370
$ws.Range("J107").Value = 4146.5
$ws.Range("K107").Value = 4110
$ws.Range("L107").Value = 12439.5
$ws.Range("M107").Value = -2190
$ws.Range("N107").Value = -16279.5
$ws.Range("H122").Value = 3765.8333
$ws.Range("I122").Value = 2700.3333
$ws.Range("J122").Value = 5541.6665
$ws.Range("K122").Value = 8100.999899999999
$ws.Range("L122").Value = 16624.9995
$ws.Range("M122").Value = -5650.999899999999
$ws.Range("N122").Value = -21524.9995
$ws.Range("H136").Value = 3743.2415
$ws.Range("I136").Value = 2763.375
$ws.Range("J136").Value = 8446.6
$ws.Range("K136").Value = 8290.125
$ws.Range("L136").Value = 25339.8
$ws.Range("M136").Value = -5740.125
$ws.Range("N136").Value = -30439.8
